$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("April")

# Clear the old row 17 data
$ws.Range("A17:G17").ClearContents()

# Write the new row 18 data
$ws.Range("A18").Value = "17.04.2018, Tue"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
